# Applies the crypto price/volume refresh + row reshuffle described in the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: cell address -> new value (as it should literally read in the sheet).
$updates = [ordered]@{
    'D2' = '308.09'
    'E2' = '1.53%'
    'D3' = '39.10'
    'E3' = '9.49%'
    'E4' = '1.19%'
    'D5' = '0.08155'
    'E5' = '3.29%'
    'D6' = '1.986'
    'E6' = '7.50%'
    'B7' = 'GateToken'
    'C7' = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
    'D7' = '4.175'
    'E7' = '1.67%'
    'B8' = 'KuCoinToken'
    'C8' = 'https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs'
    'D8' = '7.901'
    'E8' = '1.49%'
    'B9' = 'MXToken'
    'C9' = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
    'D9' = '0.9295'
    'E9' = '1.10%'
    'B10' = 'LiechtensteinCryptoassetsExchange'
    'C10' = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
    'D10' = '0.1415'
    'E10' = '5.27%'
    'B11' = 'WazirX'
    'C11' = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
    'D11' = '0.1940'
    'E11' = '2.37%'
    'B12' = 'MandalaExchangeToken'
    'C12' = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
    'D12' = '0.09240'
    'E12' = '1.47%'
    'B13' = 'BitrueCoin'
    'C13' = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
    'D13' = '0.03481'
    'E13' = '0.30%'
    'B14' = 'BitMartToken'
    'C14' = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
    'D14' = '0.09831'
    'E14' = '0.09%'
    'B15' = 'BitForexToken'
    'C15' = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
    'D15' = '0.001406'
    'E15' = '0.19%'
    'B16' = 'TigerCash'
    'C16' = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
    'D16' = '0.005950'
    'E16' = '-2.34%'
    'B17' = 'LEO'
    'C17' = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
    'D17' = '3.947'
    'E17' = '6.00%'
    'E18' = '2.66%'
    'D19' = '0.3453'
    'E19' = '0.44%'
    'E20' = '-0.57%'
    'D21' = '4.800'
    'E22' = '19.38%'
    'D23' = '0.04472'
    'E23' = '1.52%'
    'D24' = '0.001243'
    'E24' = '0.56%'
    'E25' = '-9.58%'
    'D27' = '0.0001300'
    'E27' = '-0.10%'
    'D39' = '0.02114'
    'E39' = '9.00%'
    'D40' = '0.05164'
    'E40' = '1.85%'
    'D41' = '0.007453'
    'E41' = '-2.33%'
    'D42' = '0.01014'
    'E42' = '-0.56%'
    'D43' = '0.1367'
    'E43' = '1.95%'
    'D44' = '0.002130'
    'E44' = '-1.49%'
    'D45' = '0.009668'
    'E45' = '-5.09%'
    'D46' = '0.00006302'
    'E46' = '2.24%'
    'D47' = '0.00000000750'
    'E47' = '-0.11%'
    'E48' = '1.94%'
    'D49' = '0.001600'
    'E49' = '-3.65%'
    'D50' = '0.00002100'
    'E50' = '-0.11%'
    'D51' = '0.0002000'
    'E51' = '-0.11%'
}

# Columns D (Price) and E (Volume(1h)) hold plain text like "39.10" or "1.53%".
# Force text format first so Excel does not silently reinterpret them as numbers
# (which would drop trailing zeros and turn percentages into decimal fractions).
foreach ($addr in $updates.Keys) {
    $col = ($addr -replace '[0-9]+$', '')
    $cell = $ws.Range($addr)
    if ($col -eq 'D' -or $col -eq 'E') {
        $cell.NumberFormat = '@'
    }
    $cell.Value = $updates[$addr]
}
